# results_2hl.xlsx update:
#  - rows 2-25 (existing runs): re-randomized "init" folds -> new D (init col),
#    C (num_nodes, only for rows 6-25) and recomputed H (mse) / I (r2) values
#  - rows 26-37: twelve new runs appended (num_nodes 200 and 400, each with
#    the same six "init" values used throughout the sheet)
#  - selection moves to the last-entered cell, K36

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- updates to existing rows 2-25 -----------------------------------------
# Row, num_nodes(C), init(D), mse(H), r2(I)
# NOTE: each inner array is prefixed with the unary comma operator so the
# outer @() keeps one element per row instead of flattening everything into
# a single flat list.
$updates = @(
    ,@(2,   20, 3,  342.44914799999998,  0.79490499999999997)
    ,@(3,   20, 5,  241.04107099999999,  0.85563900000000004)
    ,@(4,   20, 6,  357.04757999999998,  0.78616200000000003)
    ,@(5,   20, 8,  331.01748900000001,  0.80175200000000002)
    ,@(6,   20, 9,  35.199548999999998,  0.97891899999999998)
    ,@(7,   20, 20, 83.215029999999999,  0.95016199999999995)
    ,@(8,   24, 3,  233.03638100000001,  0.860433)
    ,@(9,   24, 5,  192.958068,          0.884436)
    ,@(10,  24, 6,  291.492706,          0.82542300000000002)
    ,@(11,  24, 8,  74.974761999999998,  0.95509699999999997)
    ,@(12,  24, 9,  52.05527,            0.96882400000000002)
    ,@(13,  24, 20, 336.75471499999998,  0.79831600000000003)
    ,@(14,  50, 3,  343.05552399999999,  0.79454199999999997)
    ,@(15,  50, 5,  431.70894700000002,  0.74144699999999997)
    ,@(16,  50, 6,  384.33625799999999,  0.76981900000000003)
    ,@(17,  50, 8,  484.59704199999999,  0.70977199999999996)
    ,@(18,  50, 9,  587.20713699999999,  0.64831799999999995)
    ,@(19,  50, 20, 225.75714099999999,  0.86479300000000003)
    ,@(20,  100, 3, 4460.1546609999996, -1.671211)
    ,@(21,  100, 5, 315.17017499999997,  0.81124300000000005)
    ,@(22,  100, 6, 262.19066800000002,  0.84297299999999997)
    ,@(23,  100, 8, 193.92537300000001,  0.883857)
    ,@(24,  100, 9, 338.59036500000002,  0.79721600000000004)
    ,@(25,  100, 20, 225.739735,         0.86480299999999999)
)

foreach ($u in $updates) {
    $r = $u[0]
    $ws.Cells.Item($r, 3).Value = $u[1]   # C: num_nodes
    $ws.Cells.Item($r, 4).Value = $u[2]   # D: init
    $ws.Cells.Item($r, 8).Value = $u[3]   # H: mse
    $ws.Cells.Item($r, 9).Value = $u[4]   # I: r2
}

# --- twelve new rows (26-37) ------------------------------------------------
# Row, run(B), num_nodes(C), init(D), mse(H), r2(I)
$newRows = @(
    ,@(26, 25, 200, 3,  317.27184899999997,  0.80998400000000004)
    ,@(27, 26, 200, 5,  315.92709200000002,  0.81079000000000001)
    ,@(28, 27, 200, 6,  359.04607700000003,  0.78496500000000002)
    ,@(29, 28, 200, 8,  369.34921700000001,  0.77879500000000002)
    ,@(30, 29, 200, 9,  241.43160599999999,  0.85540499999999997)
    ,@(31, 30, 200, 20, 278.20194099999998,  0.83338299999999998)
    ,@(32, 31, 400, 3,  4585.9464239999998, -1.746548)
    ,@(33, 32, 400, 5,  613.44965500000001,  0.632602)
    ,@(34, 33, 400, 6,  564.929709,          0.66166100000000005)
    ,@(35, 34, 400, 8,  532.91030799999999,  0.68083700000000003)
    ,@(36, 35, 400, 9,  502.42801500000002,  0.69909299999999996)
    ,@(37, 36, 400, 20, 410.73434200000003,  0.75400900000000004)
)

foreach ($n in $newRows) {
    $r = $n[0]
    $ws.Cells.Item($r, 1).Value = 0          # A: run
    $ws.Cells.Item($r, 2).Value = $n[1]      # B: k
    $ws.Cells.Item($r, 3).Value = $n[2]      # C: num_nodes
    $ws.Cells.Item($r, 4).Value = $n[3]      # D: init
    $ws.Cells.Item($r, 5).Value = "relu"     # E: activation_func
    $ws.Cells.Item($r, 6).Value = 0          # F: dropout_rate
    $ws.Cells.Item($r, 7).Value = "normal"   # G: init (weight init)
    $ws.Cells.Item($r, 8).Value = $n[4]      # H: mse
    $ws.Cells.Item($r, 9).Value = $n[5]      # I: r2
}

# --- match the saved selection in the source workbook ----------------------
$ws.Range("K36").Select() | Out-Null

Write-Output "Updated rows 2-25 and appended rows 26-37 on $($ws.Name)."
